$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 8272.091
$ws.Range("I2").Value = 7285.4287
$ws.Range("J2").Value = 9998.75
$ws.Range("K2").Value = 7285.4287
$ws.Range("L2").Value = 9998.75
$ws.Range("M2").Value = -7172.4287
$ws.Range("N2").Value = -10224.75
$ws.Range("H5").Value = 245
$ws.Range("I5").Value = 278
$ws.Range("K5").Value = 278
$ws.Range("M5").Value = -163
$ws.Range("H34").Value = 2146.8
$ws.Range("I34").Value = 2146.8
$ws.Range("K34").Value = 2146.8
$ws.Range("M34").Value = -1943.8
$ws.Range("H36").Value = 2146.8
$ws.Range("I36").Value = 2146.8
$ws.Range("K36").Value = 2146.8
$ws.Range("M36").Value = -1431.8
$ws.Range("H40").Value = 1143.6666
$ws.Range("I40").Value = 721
$ws.Range("K40").Value = 721
$ws.Range("M40").Value = -546
$ws.Range("H98").Value = 3245
$ws.Range("I98").Value = 2101.8635
$ws.Range("K98").Value = 2101.8635
$ws.Range("M98").Value = -603.8634999999999
$ws.Range("H113").Value = 3953
$ws.Range("I113").Value = 3940.4285
$ws.Range("K113").Value = 3940.4285
$ws.Range("M113").Value = -686.4285
$ws.Range("H122").Value = 3245
$ws.Range("I122").Value = 2101.8635
$ws.Range("K122").Value = 6305.5905
$ws.Range("M122").Value = -3855.5905
$ws.Range("H137").Value = 25000538
$ws.Range("I137").Value = 275
$ws.Range("J137").Value = 50000800
$ws.Range("K137").Value = 825
$ws.Range("L137").Value = 150002400
$ws.Range("M137").Value = 1725
$ws.Range("N137").Value = -150007500

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 165453.33
$ws.Range("I32").Value = 170468.47
$ws.Range("J32").Value = 17506.5
$ws.Range("K32").Value = 170468.47
$ws.Range("L32").Value = 17506.5
$ws.Range("M32").Value = -170181.47
$ws.Range("N32").Value = -18080.5
$ws.Range("H45").Value = 47592.137
$ws.Range("I45").Value = 60784.47
$ws.Range("K45").Value = 60784.47
$ws.Range("M45").Value = -60407.47
$ws.Range("H88").Value = 1754.1177
$ws.Range("I88").Value = 1809.75
$ws.Range("J88").Value = 1704.6666
$ws.Range("K88").Value = 1809.75
$ws.Range("L88").Value = 1704.6666
$ws.Range("M88").Value = -1403.75
$ws.Range("N88").Value = -2516.6666
$ws.Range("H91").Value = 1754.1177
$ws.Range("I91").Value = 1809.75
$ws.Range("J91").Value = 1704.6666
$ws.Range("K91").Value = 1809.75
$ws.Range("L91").Value = 1704.6666
$ws.Range("M91").Value = -405.75
$ws.Range("N91").Value = -4512.6666
$ws.Range("H122").Value = 1072.4286
$ws.Range("I122").Value = 1072.4286
$ws.Range("K122").Value = 3217.2858
$ws.Range("M122").Value = -767.2857999999997
$ws.Range("H132").Value = 2137.2449
$ws.Range("I132").Value = 2015.4255
$ws.Range("K132").Value = 6046.2765
$ws.Range("M132").Value = -3516.2765

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1387.9231
$ws.Range("I22").Value = 1094.8182
$ws.Range("K22").Value = 1094.8182
$ws.Range("M22").Value = -921.8181999999999
$ws.Range("H134").Value = 23686448
$ws.Range("I134").Value = 1994.3438
$ws.Range("K134").Value = 5983.0314
$ws.Range("M134").Value = -3448.0314

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 11050
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -28744

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 45009140
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 45009140
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = $null
$ws.Range("M9").Value = 135027420
$ws.Range("N9").Value = -135027868
$ws.Range("H33").Value = 25001548
$ws.Range("I33").Value = 600.75
$ws.Range("J33").Value = 50002496
$ws.Range("K33").Value = 3604.5
$ws.Range("L33").Value = 300014976
$ws.Range("M33").Value = -3321.5
$ws.Range("N33").Value = -300015542
$ws.Range("H37").Value = 124995
$ws.Range("J37").Value = 124995
$ws.Range("L37").Value = 374985
$ws.Range("N37").Value = -375209
$ws.Range("H68").Value = 2531.375
$ws.Range("J68").Value = 2902
$ws.Range("L68").Value = 8706
$ws.Range("N68").Value = -10328
$ws.Range("H69").Value = 8394.15
$ws.Range("J69").Value = 6493.3125
$ws.Range("L69").Value = 19479.9375
$ws.Range("N69").Value = -21101.9375
$ws.Range("H71").Value = 2531.375
$ws.Range("J71").Value = 2902
$ws.Range("L71").Value = 26118
$ws.Range("N71").Value = -34230
$ws.Range("H72").Value = 8394.15
$ws.Range("J72").Value = 6493.3125
$ws.Range("L72").Value = 58439.8125
$ws.Range("N72").Value = -66551.8125
$ws.Range("H92").Value = 519.8
$ws.Range("I92").Value = 449.75
$ws.Range("K92").Value = 1349.25
$ws.Range("M92").Value = -101.25
$ws.Range("H131").Value = 55640.527
$ws.Range("I131").Value = 3425.3333
$ws.Range("J131").Value = 79739.84
$ws.Range("K131").Value = 10275.9999
$ws.Range("L131").Value = 239219.52
$ws.Range("M131").Value = -5235.999899999999
$ws.Range("N131").Value = -249299.52

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2749.0833
$ws.Range("I102").Value = 2428.9
$ws.Range("J102").Value = 4350
$ws.Range("K102").Value = 2428.9
$ws.Range("L102").Value = 4350
$ws.Range("M102").Value = -806.9000000000001
$ws.Range("N102").Value = -7594
$ws.Range("H122").Value = 2761.5483
$ws.Range("I122").Value = 2683.7083
$ws.Range("K122").Value = 8051.124899999999
$ws.Range("M122").Value = -5601.124899999999
$ws.Range("H126").Value = 2197.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7096.3
$ws.Range("I7").Value = 3381.7693
$ws.Range("K7").Value = 3381.7693
$ws.Range("M7").Value = -3269.7693
$ws.Range("H40").Value = 2981.0588
$ws.Range("I40").Value = 1898.909
$ws.Range("K40").Value = 1898.909
$ws.Range("M40").Value = -1762.909
$ws.Range("H43").Value = 2598125
$ws.Range("J43").Value = 4130000
$ws.Range("L43").Value = 4130000
$ws.Range("N43").Value = -4130386
$ws.Range("H46").Value = 14272.75
$ws.Range("I46").Value = 37291.332
$ws.Range("J46").Value = 6599.8887
$ws.Range("K46").Value = 37291.332
$ws.Range("L46").Value = 6599.8887
$ws.Range("M46").Value = -37103.332
$ws.Range("N46").Value = -6975.8887
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = $null
$ws.Range("N119").Value = 0
$ws.Range("H122").Value = 3383
$ws.Range("J122").Value = 4153.8125
$ws.Range("L122").Value = 12461.4375
$ws.Range("N122").Value = -17361.4375
$ws.Range("H126").Value = 7096.3
$ws.Range("I126").Value = 3381.7693
$ws.Range("K126").Value = 10145.3079
$ws.Range("M126").Value = -7675.3079
$ws.Range("H132").Value = 3128.9707
$ws.Range("I132").Value = 2755.7273
$ws.Range("K132").Value = 8267.1819
$ws.Range("M132").Value = -5737.1819
$ws.Range("H136").Value = 2270.4143
$ws.Range("I136").Value = 3798.4285
$ws.Range("J136").Value = 1888.4108
$ws.Range("K136").Value = 11395.2855
$ws.Range("L136").Value = 5665.232400000001
$ws.Range("M136").Value = -8845.2855
$ws.Range("N136").Value = -10765.2324

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = $null
$ws.Range("H30").Value = 4997.5
$ws.Range("J30").Value = 4997.5
$ws.Range("L30").Value = 4997.5
$ws.Range("N30").Value = -5211.5
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 10000
$ws.Range("K39").Value = 10000
$ws.Range("M39").Value = -9587
$ws.Range("H62").Value = 2126.9333
$ws.Range("I62").Value = 2299.182
$ws.Range("J62").Value = 1653.25
$ws.Range("K62").Value = 2299.182
$ws.Range("L62").Value = 1653.25
$ws.Range("M62").Value = -1675.182
$ws.Range("N62").Value = -2901.25
$ws.Range("H65").Value = 2126.9333
$ws.Range("I65").Value = 2299.182
$ws.Range("J65").Value = 1653.25
$ws.Range("K65").Value = 11495.91
$ws.Range("L65").Value = 8266.25
$ws.Range("M65").Value = -8375.91
$ws.Range("N65").Value = -14506.25
$ws.Range("H74").Value = 12799
$ws.Range("J74").Value = 23500
$ws.Range("L74").Value = 23500
$ws.Range("N74").Value = -25372
$ws.Range("H77").Value = 12799
$ws.Range("J77").Value = 23500
$ws.Range("L77").Value = 70500
$ws.Range("N77").Value = -79860
$ws.Range("H107").Value = 893450.3
$ws.Range("I107").Value = 608.0357
$ws.Range("K107").Value = 1824.1071
$ws.Range("M107").Value = 95.89289999999983
$ws.Range("H122").Value = 4600
$ws.Range("I122").Value = 3250
$ws.Range("K122").Value = 9750
$ws.Range("M122").Value = -7300
$ws.Range("H126").Value = 2979.3572
$ws.Range("I126").Value = 2785.3
$ws.Range("J126").Value = 3464.5
$ws.Range("K126").Value = 8355.900000000001
$ws.Range("L126").Value = 10393.5
$ws.Range("M126").Value = -5885.900000000001
$ws.Range("N126").Value = -15333.5

Write-Host "Applied all changes"